# Updated symbol list on Sun Dec 11 23:36:20 UTC 2022 with GitHub Actions
#
# Column D ("Price") holds numeric-looking values that are stored as TEXT
# (inline strings) in the workbook, not as numbers. Assigning a plain
# numeric-looking string to a General-formatted cell makes Excel coerce it
# to a real number (and can introduce floating point noise, e.g.
# 284.97 -> 284.97000000000003), which would not match the source data.
#
# To keep these as text we prefix the literal with an apostrophe (Excel's
# "treat as text" quote-prefix). That correctly keeps the value as a
# string, but it also marks the cell's style with quotePrefix="1" (visible
# to Excel as a small green corner marker). Since the original cells carry
# no such marker, we immediately copy the number format/style back from an
# untouched neighboring cell in the same column so the cell's style index
# reverts to its original (unflagged) state while the text content stays.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($row, $value) {
    $cell = $ws.Cells.Item($row, 4)   # column D
    $cell.Value = "'" + $value
    # Restore the original (unflagged) style by copying from a nearby,
    # unmodified cell in the same column.
    if ($row -eq 2) {
        $ref = $ws.Cells.Item(3, 4)
    } else {
        $ref = $ws.Cells.Item($row - 1, 4)
    }
    $cell.Style = $ref.Style
}

Set-TextPrice 2  "284.97"
Set-TextPrice 3  "21.20"
Set-TextPrice 4  "6.451"
Set-TextPrice 5  "0.06356"
Set-TextPrice 6  "3.600"
Set-TextPrice 7  "1.537"
Set-TextPrice 9  "0.8206"
Set-TextPrice 10 "0.01411"
Set-TextPrice 11 "0.1672"
Set-TextPrice 12 "0.08679"
Set-TextPrice 13 "0.03669"
Set-TextPrice 14 "0.03229"
Set-TextPrice 15 "0.09197"
Set-TextPrice 16 "3.706"
Set-TextPrice 17 "0.001649"
Set-TextPrice 18 "0.04737"
Set-TextPrice 19 "0.006249"
Set-TextPrice 22 "0.0001603"
Set-TextPrice 23 "3.785"
Set-TextPrice 24 "2.265"
Set-TextPrice 25 "0.3356"
Set-TextPrice 40 "0.04761"

# Row 41: KickToken -> BKEXToken (values swapped with row 43 below)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice 41 "0.1117"
$ws.Range("E41").Value = "40BKEXTokenBKK"

Set-TextPrice 42 "0.003455"

# Row 43: BKEXToken -> KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice 43 "0.003584"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

Set-TextPrice 44 "0.01169"
Set-TextPrice 45 "0.00006935"
Set-TextPrice 47 "1.001"
Set-TextPrice 48 "0.004187"

$ws.Range("E49").Value = "48CryptobidCoinCBC"
